# Append the latest Adafruit IO reading as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Value") must stay text (matches the existing text-typed rows),
# so force text formatting before writing a numeric-looking string into it.
$ws.Range("C47").NumberFormat = "@"

$ws.Range("A47").Value = "2024-09-25T18:06:40Z"
$ws.Range("B47").Value = "temperature"
$ws.Range("C47").Value = "25"
$ws.Range("D47").Value = "N/A"
$ws.Range("E47").Value = "N/A"
$ws.Range("F47").Value = "N/A"
